$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently holds subjects S1-S4 in rows 2-61 (row 1 is the header).
# This commit appends a 5th subject (S5 / Jaime, session "2013-03-05-jaime")
# as 15 new data rows (62-76), one per recording file, growing the used range
# from A1:G61 to A1:G76.

# Seed the new rows with the same per-row formatting used by the existing data
# (centered text, and a short-date format on column C) by copying row 61's
# formats down across the new block, then fill in the values below.
$ws.Range("A61:G61").Copy()
$ws.Range("A62:G76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column E (fileName) is explicitly centered in the target workbook; make sure
# that holds even though the paste above already centers it.
$ws.Range("E62:E76").HorizontalAlignment = -4108

$ws.Range("A62").Value2 = "S5"
$ws.Range("B62").Value2 = "Jaime"
$ws.Range("C62").Value2 = 41338
$ws.Range("D62").Value2 = "2013-03-05-jaime"
$ws.Range("E62").Value2 = "2013-03-05-16-22-01"
$ws.Range("F62").Value2 = "hybrid-15Hz"
$ws.Range("G62").Value2 = 1
$ws.Range("A63").Value2 = "S5"
$ws.Range("B63").Value2 = "Jaime"
$ws.Range("C63").Value2 = 41338
$ws.Range("D63").Value2 = "2013-03-05-jaime"
$ws.Range("E63").Value2 = "2013-03-05-16-30-48"
$ws.Range("F63").Value2 = "hybrid-8-57Hz"
$ws.Range("G63").Value2 = 1
$ws.Range("A64").Value2 = "S5"
$ws.Range("B64").Value2 = "Jaime"
$ws.Range("C64").Value2 = 41338
$ws.Range("D64").Value2 = "2013-03-05-jaime"
$ws.Range("E64").Value2 = "2013-03-05-16-38-38"
$ws.Range("F64").Value2 = "hybrid-15Hz"
$ws.Range("G64").Value2 = 2
$ws.Range("A65").Value2 = "S5"
$ws.Range("B65").Value2 = "Jaime"
$ws.Range("C65").Value2 = 41338
$ws.Range("D65").Value2 = "2013-03-05-jaime"
$ws.Range("E65").Value2 = "2013-03-05-16-44-25"
$ws.Range("F65").Value2 = "oddball.bdf"
$ws.Range("G65").Value2 = 1
$ws.Range("A66").Value2 = "S5"
$ws.Range("B66").Value2 = "Jaime"
$ws.Range("C66").Value2 = 41338
$ws.Range("D66").Value2 = "2013-03-05-jaime"
$ws.Range("E66").Value2 = "2013-03-05-16-51-03"
$ws.Range("F66").Value2 = "hybrid-12Hz"
$ws.Range("G66").Value2 = 1
$ws.Range("A67").Value2 = "S5"
$ws.Range("B67").Value2 = "Jaime"
$ws.Range("C67").Value2 = 41338
$ws.Range("D67").Value2 = "2013-03-05-jaime"
$ws.Range("E67").Value2 = "2013-03-05-17-04-04"
$ws.Range("F67").Value2 = "hybrid-10Hz"
$ws.Range("G67").Value2 = 1
$ws.Range("A68").Value2 = "S5"
$ws.Range("B68").Value2 = "Jaime"
$ws.Range("C68").Value2 = 41338
$ws.Range("D68").Value2 = "2013-03-05-jaime"
$ws.Range("E68").Value2 = "2013-03-05-17-09-20"
$ws.Range("F68").Value2 = "hybrid-10Hz"
$ws.Range("G68").Value2 = 2
$ws.Range("A69").Value2 = "S5"
$ws.Range("B69").Value2 = "Jaime"
$ws.Range("C69").Value2 = 41338
$ws.Range("D69").Value2 = "2013-03-05-jaime"
$ws.Range("E69").Value2 = "2013-03-05-17-17-35"
$ws.Range("F69").Value2 = "oddball.bdf"
$ws.Range("G69").Value2 = 2
$ws.Range("A70").Value2 = "S5"
$ws.Range("B70").Value2 = "Jaime"
$ws.Range("C70").Value2 = 41338
$ws.Range("D70").Value2 = "2013-03-05-jaime"
$ws.Range("E70").Value2 = "2013-03-05-17-23-14"
$ws.Range("F70").Value2 = "hybrid-12Hz"
$ws.Range("G70").Value2 = 2
$ws.Range("A71").Value2 = "S5"
$ws.Range("B71").Value2 = "Jaime"
$ws.Range("C71").Value2 = 41338
$ws.Range("D71").Value2 = "2013-03-05-jaime"
$ws.Range("E71").Value2 = "2013-03-05-17-38-14"
$ws.Range("F71").Value2 = "hybrid-15Hz"
$ws.Range("G71").Value2 = 3
$ws.Range("A72").Value2 = "S5"
$ws.Range("B72").Value2 = "Jaime"
$ws.Range("C72").Value2 = 41338
$ws.Range("D72").Value2 = "2013-03-05-jaime"
$ws.Range("E72").Value2 = "2013-03-05-17-43-55"
$ws.Range("F72").Value2 = "oddball.bdf"
$ws.Range("G72").Value2 = 3
$ws.Range("A73").Value2 = "S5"
$ws.Range("B73").Value2 = "Jaime"
$ws.Range("C73").Value2 = 41338
$ws.Range("D73").Value2 = "2013-03-05-jaime"
$ws.Range("E73").Value2 = "2013-03-05-17-49-15"
$ws.Range("F73").Value2 = "hybrid-10Hz"
$ws.Range("G73").Value2 = 3
$ws.Range("A74").Value2 = "S5"
$ws.Range("B74").Value2 = "Jaime"
$ws.Range("C74").Value2 = 41338
$ws.Range("D74").Value2 = "2013-03-05-jaime"
$ws.Range("E74").Value2 = "2013-03-05-17-55-29"
$ws.Range("F74").Value2 = "hybrid-8-57Hz"
$ws.Range("G74").Value2 = 2
$ws.Range("A75").Value2 = "S5"
$ws.Range("B75").Value2 = "Jaime"
$ws.Range("C75").Value2 = 41338
$ws.Range("D75").Value2 = "2013-03-05-jaime"
$ws.Range("E75").Value2 = "2013-03-05-18-00-47"
$ws.Range("F75").Value2 = "hybrid-8-57Hz"
$ws.Range("G75").Value2 = 3
$ws.Range("A76").Value2 = "S5"
$ws.Range("B76").Value2 = "Jaime"
$ws.Range("C76").Value2 = 41338
$ws.Range("D76").Value2 = "2013-03-05-jaime"
$ws.Range("E76").Value2 = "2013-03-05-18-06-33"
$ws.Range("F76").Value2 = "hybrid-12Hz"
$ws.Range("G76").Value2 = 3

# Match the post-edit cursor/selection position recorded in the workbook.
$ws.Range("A77").Select()
